$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 21 values (I21, J21)
$ws.Range("I21").Value = 192
$ws.Range("J21").Value = 608

# Add new row 40 - "Frozen Soul" NPC, promotional material for the corrupted bishop
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "FrozenSoul"
$ws.Range("C40").Value = "Frozen Soul"
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = "The Ice Plane"
$ws.Range("I40").Value = 336
$ws.Range("J40").Value = 1248
